$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.166.84"
$ws.Range("E2").Value = "  +1.38%  "

$ws.Range("D3").Value = "2.596.33"
$ws.Range("E3").Value = "  +0.33%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'529.94"
$ws.Range("E5").Value = "  +2.17%  "

$ws.Range("D6").Value = "'140.44"
$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("E8").Value = "  +0.62%  "

$ws.Range("D9").Value = "2.610.03"
$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("D10").Value = "'6.45"
$ws.Range("E10").Value = "  -0.21%  "

$ws.Range("E11").Value = "  +1.09%  "

$ws.Range("E12").Value = "  +0.86%  "

$ws.Range("E13").Value = "  +2.97%  "

$ws.Range("D14").Value = "3.056.51"
$ws.Range("E14").Value = "  +0.60%  "

$ws.Range("D15").Value = "59.113.24"
$ws.Range("E15").Value = "  +1.44%  "

$ws.Range("D16").Value = "'20.52"
$ws.Range("E16").Value = "  +1.01%  "

$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("D18").Value = "2.595.74"
$ws.Range("E18").Value = "  +0.41%  "

$ws.Range("D19").Value = "'347.45"
$ws.Range("E19").Value = "  +2.92%  "

$ws.Range("D20").Value = "'4.34"
$ws.Range("E20").Value = "  +0.97%  "

$ws.Range("E21").Value = "  -0.79%  "

$ws.Range("E22").Value = "  +0.58%  "

$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'67.36"
$ws.Range("E24").Value = "  +3.04%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("E26").Value = "  +1.34%  "

$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("E28").Value = "  +2.31%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0741"
$ws.Range("E29").Value = "  +0.99%  "

$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("E31").Value = "  +3.26%  "

$ws.Range("D32").Value = "'5.89"
$ws.Range("E32").Value = "  -3.67%  "

$ws.Range("D33").Value = "'18.80"
$ws.Range("E33").Value = "  +0.71%  "

$ws.Range("D34").Value = "'149.21"
$ws.Range("E34").Value = "  -0.24%  "

$ws.Range("E35").Value = "  +0.63%  "

$ws.Range("D36").Value = "'1.13"
$ws.Range("E36").Value = "  -0.28%  "

$ws.Range("D37").Value = "'36.81"
$ws.Range("E37").Value = "  +1.85%  "

$ws.Range("D38").Value = "'1.47"
$ws.Range("E38").Value = "  +0.44%  "

$ws.Range("E39").Value = "  -2.61%  "

$ws.Range("D40").Value = "'0.822"
$ws.Range("E40").Value = "  -0.69%  "

$ws.Range("E41").Value = "  +1.18%  "

$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").Value = "'0.598"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").Value = "'269.73"
$ws.Range("E44").Value = "  -1.77%  "

$ws.Range("D45").Value = "'10.77"
$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("D46").Value = "'0.0959"
$ws.Range("E46").Value = "  +1.37%  "

$ws.Range("D47").Value = "'0.0520"
$ws.Range("E47").Value = "  +0.05%  "

$ws.Range("D48").Value = "'18.47"
$ws.Range("E48").Value = "  -1.41%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.63"
$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.963.32"
$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0222"
$ws.Range("E51").Value = "  +0.89%  "
